$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Resolve the product/supplier reference number (A2). It is a numeric-looking
# string ("03250001") that must be kept as text (leading zero), so force a
# text format while entering it, then restore the cell style so no extra
# formatting is left behind.
$cellA2 = $ws.Range("A2")
$cellA2.NumberFormat = "@"
$cellA2.Value = "03250001"
$cellA2.Style = "Normal"

# Update the creation date/time for the remaining incident row.
$ws.Range("B2").Value = 45720.35496090278

# Fill in the previously empty description and update the status.
$ws.Range("I2").Value = "dfgdf"
$ws.Range("K2").Value = "EN ATTENTE"

# The second incident row is no longer needed; remove it entirely.
$ws.Rows("3:3").Delete()
